$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("B2").Value = "H2JKV"
$ws.Range("C2").Value = 12621

# Remove row 3 entirely
$ws.Rows.Item(3).Delete()
